# Feature - Feedback button for retraining
# Applies corrected emotion/sub_emotion labels, a couple of text fixes, and
# splits one mis-merged transcript line into two rows (shifting the tail of
# the sheet down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple in-place label corrections (no row shift needed) ---

# Row 19
$ws.Range("E19").Value = "admiration"

# Row 82
$ws.Range("D82").Value = "fear"
$ws.Range("E82").Value = "nervousness"

# Row 89
$ws.Range("E89").Value = "curiosity"

# Row 111
$ws.Range("D111").Value = "fear"
$ws.Range("E111").Value = "nervousness"

# Row 120
$ws.Range("E120").Value = "excitement"

# Row 144
$ws.Range("C144").Value = "Watching."
$ws.Range("D144").Value = "neutral"
$ws.Range("E144").Value = "neutral"

# Row 163
$ws.Range("D163").Value = "fear"
$ws.Range("E163").Value = "nervousness"

# Row 164 - also fix the end-time and split the merged sentence
$ws.Range("B164").Value = "00:13:46"
$ws.Range("C164").Value = "Can you fly back?"
$ws.Range("D164").Value = "fear"
$ws.Range("E164").Value = "nervousness"

# --- Insert a new row 165 for the split-off second half of the sentence ---
# This shifts the previous rows 165-182 down to 166-183.
$ws.Rows.Item(165).Insert()

$ws.Range("A165").Value = "00:13:46"
$ws.Range("B165").Value = "00:13:47"
$ws.Range("C165").Value = "He fly up?"
$ws.Range("D165").Value = "fear"
$ws.Range("E165").Value = "nervousness"
$ws.Range("F165").Value = "mild"

# --- Fix a capitalization typo that survived the shift (old row180 -> new row181) ---
$ws.Range("C181").Value = "Special thanks to everyone at Kansat who made this experience possible."
